# chore: update Sheets via scheduled runner
# Refresh cached market-board price/profit figures on the per-job "Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2137.1
$ws.Range("J40").Value = 2361.8333
$ws.Range("L40").Value = 2361.8333
$ws.Range("N40").Value = -2711.8333

$ws.Range("H69").Value = 191441
$ws.Range("I69").Value = 19199.8
$ws.Range("J69").Value = 263208.16
$ws.Range("K69").Value = 57599.39999999999
$ws.Range("L69").Value = 789624.48
$ws.Range("M69").Value = -56725.39999999999
$ws.Range("N69").Value = -791372.48

$ws.Range("H72").Value = 191441
$ws.Range("I72").Value = 19199.8
$ws.Range("J72").Value = 263208.16
$ws.Range("K72").Value = 172798.2
$ws.Range("L72").Value = 2368873.44
$ws.Range("M72").Value = -168430.2
$ws.Range("N72").Value = -2377609.44

$ws.Range("H80").Value = 4771
$ws.Range("J80").Value = 6959.778
$ws.Range("L80").Value = 20879.334
$ws.Range("N80").Value = -22875.334

$ws.Range("H83").Value = 4771
$ws.Range("J83").Value = 6959.778
$ws.Range("L83").Value = 62638.002
$ws.Range("N83").Value = -72622.00200000001

$ws.Range("H137").Value = 22622.941
$ws.Range("I137").Value = 103000.336
$ws.Range("J137").Value = 14844.483
$ws.Range("K137").Value = 309001.008
$ws.Range("L137").Value = 44533.449
$ws.Range("M137").Value = -306451.008
$ws.Range("N137").Value = -49633.449

$ws.Range("H141").Value = 1461.1765
$ws.Range("I141").Value = 1427.5
$ws.Range("K141").Value = 4282.5
$ws.Range("M141").Value = 897.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4842
$ws.Range("I45").Value = 3999
$ws.Range("J45").Value = 5071.909
$ws.Range("K45").Value = 3999
$ws.Range("L45").Value = 5071.909
$ws.Range("M45").Value = -3622
$ws.Range("N45").Value = -5825.909

$ws.Range("H61").Value = 7127.125
$ws.Range("I61").Value = 841
$ws.Range("K61").Value = 841
$ws.Range("M61").Value = -629

$ws.Range("H74").Value = 764676.1
$ws.Range("I74").Value = 1201481.8
$ws.Range("K74").Value = 1201481.8
$ws.Range("M74").Value = -1200607.8

$ws.Range("H77").Value = 764676.1
$ws.Range("I77").Value = 1201481.8
$ws.Range("K77").Value = 6007409
$ws.Range("M77").Value = -6003041

$ws.Range("H97").Value = 639.7059
$ws.Range("I97").Value = 577.62067
$ws.Range("K97").Value = 577.62067
$ws.Range("M97").Value = -81.62067000000002

$ws.Range("H132").Value = 1780.4722
$ws.Range("I132").Value = 1447.9
$ws.Range("J132").Value = 3443.3333
$ws.Range("K132").Value = 4343.700000000001
$ws.Range("L132").Value = 10329.9999
$ws.Range("M132").Value = -1813.700000000001
$ws.Range("N132").Value = -15389.9999

$ws.Range("H136").Value = 7127.125
$ws.Range("I136").Value = 841
$ws.Range("K136").Value = 2523
$ws.Range("M136").Value = 27

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1088.6086
$ws.Range("I58").Value = 1004.3158
$ws.Range("K58").Value = 1004.3158
$ws.Range("M58").Value = -801.3158

$ws.Range("H99").Value = 1889.9231
$ws.Range("I99").Value = 1508.8889
$ws.Range("K99").Value = 1508.8889
$ws.Range("M99").Value = -10.88889999999992

$ws.Range("H126").Value = 1889.9231
$ws.Range("I126").Value = 1508.8889
$ws.Range("K126").Value = 4526.6667
$ws.Range("M126").Value = -2056.6667

$ws.Range("H132").Value = 67601.2
$ws.Range("I132").Value = 72358.42999999999
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 217075.29
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -214545.29
$ws.Range("N132").Value = -8060

$ws.Range("H134").Value = 3021.7334
$ws.Range("I134").Value = 2994.3076
$ws.Range("K134").Value = 8982.9228
$ws.Range("M134").Value = -6447.9228

$ws.Range("H136").Value = 1088.6086
$ws.Range("I136").Value = 1004.3158
$ws.Range("K136").Value = 3012.9474
$ws.Range("M136").Value = -462.9474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 349.6
$ws.Range("I2").Value = 412.5
$ws.Range("K2").Value = 2475
$ws.Range("M2").Value = -2362

$ws.Range("H38").Value = 33333470
$ws.Range("I38").Value = 60
$ws.Range("K38").Value = 180
$ws.Range("M38").Value = 167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7111.4736
$ws.Range("I80").Value = 2409.6365
$ws.Range("J80").Value = 13576.5
$ws.Range("K80").Value = 2409.6365
$ws.Range("L80").Value = 13576.5
$ws.Range("M80").Value = -1411.6365
$ws.Range("N80").Value = -15572.5

$ws.Range("H83").Value = 7111.4736
$ws.Range("I83").Value = 2409.6365
$ws.Range("J83").Value = 13576.5
$ws.Range("K83").Value = 12048.1825
$ws.Range("L83").Value = 67882.5
$ws.Range("M83").Value = -7056.182500000001
$ws.Range("N83").Value = -77866.5

$ws.Range("H102").Value = 26421.25
$ws.Range("I102").Value = 27758.475
$ws.Range("K102").Value = 27758.475
$ws.Range("M102").Value = -26136.475

$ws.Range("H126").Value = 2089.8108
$ws.Range("I126").Value = 1965.5294
$ws.Range("K126").Value = 5896.5882
$ws.Range("M126").Value = -3426.5882

$ws.Range("H132").Value = 3344.5
$ws.Range("I132").Value = 2316.6667
$ws.Range("J132").Value = 3785
$ws.Range("K132").Value = 6950.000100000001
$ws.Range("L132").Value = 11355
$ws.Range("M132").Value = -4420.000100000001
$ws.Range("N132").Value = -16415

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3934.5
$ws.Range("J46").Value = 7216
$ws.Range("L46").Value = 7216
$ws.Range("N46").Value = -7592

$ws.Range("H136").Value = 6788
$ws.Range("I136").Value = 8374.25
$ws.Range("J136").Value = 5730.5
$ws.Range("K136").Value = 25122.75
$ws.Range("L136").Value = 17191.5
$ws.Range("M136").Value = -22572.75
$ws.Range("N136").Value = -22291.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 39901.8
$ws.Range("I122").Value = 43515.28
$ws.Range("K122").Value = 130545.84
$ws.Range("M122").Value = -128095.84

$ws.Range("H132").Value = 56634.383
$ws.Range("I132").Value = 118975
$ws.Range("K132").Value = 356925
$ws.Range("M132").Value = -354395

$ws.Range("H136").Value = 26661.592
$ws.Range("I136").Value = 33434
$ws.Range("J136").Value = 2958.1667
$ws.Range("K136").Value = 100302
$ws.Range("L136").Value = 8874.500100000001
$ws.Range("M136").Value = -97752
$ws.Range("N136").Value = -13974.5001
